# Updated cryptos list on Sat Sep 21 16:46:42 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.128.14'
$ws.Range('E2').Value = '  +0.34%  '
$ws.Range('D3').Value = '2.558.27'
$ws.Range('E3').Value = '  +0.30%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '583.59'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.44'
$ws.Range('E6').Value = '  -0.10%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  +1.55%  '
$ws.Range('E9').Value = '  +4.07%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.64'
$ws.Range('E10').Value = '  +0.95%  '
$ws.Range('E12').Value = '  +0.87%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '27.51'
$ws.Range('E13').Value = '  +1.42%  '
$ws.Range('D14').Value = '3.015.97'
$ws.Range('E14').Value = '  +0.28%  '
$ws.Range('D15').Value = '63.039.54'
$ws.Range('E15').Value = '  +0.27%  '
$ws.Range('E16').Value = '  +5.30%  '
$ws.Range('D17').Value = '2.565.61'
$ws.Range('E17').Value = '  +1.24%  '
$ws.Range('E18').Value = '  -1.21%  '
$ws.Range('E19').Value = '  +3.85%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '341.44'
$ws.Range('E20').Value = '  +2.21%  '
$ws.Range('E21').Value = '  +0.73%  '
$ws.Range('E22').Value = '  +0.09%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '66.53'
$ws.Range('D24').Value = '2.684.43'
$ws.Range('E24').Value = '  -0.21%  '
$ws.Range('E25').Value = '  +2.62%  '
$ws.Range('E26').Value = '  +0.63%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.19'
$ws.Range('E27').Value = '  +13.34%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.52'
$ws.Range('E28').Value = '  +2.11%  '
$ws.Range('B29').Value = 'Binance-PegBSC-USD'
$ws.Range('C29').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.00'
$ws.Range('E29').Value = '  +0.17%  '
$ws.Range('B30').Value = 'SuiNetwork'
$ws.Range('C30').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.48'
$ws.Range('E30').Value = '  +0.19%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.00'
$ws.Range('E31').Value = '  +7.89%  '
$ws.Range('D32').Value = '0.0₃0826'
$ws.Range('E32').Value = '  +1.88%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '460.46'
$ws.Range('E33').Value = '  +12.10%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '176.88'
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('E35').Value = '  +2.96%  '
$ws.Range('E36').Value = '  +2.89%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.28'
$ws.Range('E37').Value = '  +2.61%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.50'
$ws.Range('E38').Value = '  +3.53%  '
$ws.Range('E39').Value = '  +0.04%  '
$ws.Range('E40').Value = '  +0.28%  '
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '150.87'
$ws.Range('E42').Value = '  -0.45%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.82'
$ws.Range('E43').Value = '  +2.36%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '21.06'
$ws.Range('E44').Value = '  +2.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0553'
$ws.Range('E45').Value = '  +6.81%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.615'
$ws.Range('E46').Value = '  +2.29%  '
$ws.Range('E47').Value = '  +2.20%  '
$ws.Range('E48').Value = '  +2.08%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '18.43'
$ws.Range('E49').Value = '  +0.45%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.74'
$ws.Range('E50').Value = '  -1.88%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '11.38'
$ws.Range('E51').Value = '  -0.19%  '
